$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.19%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.98%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.144"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.12%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07321"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.26%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.825"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'23.73%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.767"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.57%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.739"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.40%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9256"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.30%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'-1.34%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07135"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-7.27%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08010"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.72%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03006"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.13%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09919"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.44%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001488"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.06%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006239"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.77%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.457"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.77%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D19").Value = "'0.3225"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-2.28%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-1.23%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.557"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.11%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04645"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.01%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-2.64%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.24%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004733"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'6.41%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001298"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.62%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001874"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'8.55%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01715"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-4.20%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04470"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.18%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007077"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.11%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1331"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.84%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002127"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-5.46%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01044"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-21.77%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006244"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.15%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-20.76%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.920"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'171.13%"
$ws.Range("E47").Style = "Normal"
